# TestData.xlsx edit: update hrms_id test value (98 -> 75) across all sheets,
# update the occasion month (Nov -> Dec) and occasion day (30 -> 10) on the
# order pages, and update the recorded cell selections to match.

$wb = $excel.ActiveWorkbook

# --- SalesLogin ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item("SalesLogin")
$ws1.Range("A2").Value = "75"

# --- AddNewCustomerDetails ----------------------------------------------
$ws2 = $wb.Worksheets.Item("AddNewCustomerDetails")
$ws2.Range("A2").Value = "75"

# --- SearchCustomerDetails -----------------------------------------------
$ws3 = $wb.Worksheets.Item("SearchCustomerDetails")
$ws3.Range("A2").Value = "75"
$ws3.Range("A3").Value = "75"
$ws3.Range("A4").Value = "75"
$ws3.Range("A5").Value = "75"
$ws3.Range("A6").Value = "75"

# --- OrderHeaderPage -------------------------------------------------
$ws4 = $wb.Worksheets.Item("OrderHeaderPage")
$ws4.Range("A2").Value = "75"
$ws4.Range("H2").Value = "Dec"

# --- OrderDetailPage ---------------------------------------------------
$ws5 = $wb.Worksheets.Item("OrderDetailPage")
$ws5.Range("A2").Value = "75"
$ws5.Range("H2").Value = "Dec"
$ws5.Range("J2").Value = "10"

# --- Selections (match the recorded cursor positions in each sheet) ----
$ws1.Range("B2").Select()
$ws2.Range("A2").Select()
$ws3.Range("A2:A6").Select()
$ws4.Range("H2").Select()

# Leave OrderDetailPage as the active sheet/selection, as in the source file.
$ws5.Activate()
$ws5.Range("J2").Select()
